# Update the three service-port numbers in the documentation tables. Each
# number is retyped so the cell ends up holding two runs (as Word leaves
# behind when you select part of a run and type over it), and the "_GoBack"
# last-edit bookmark follows the final edit (the Cars server port cell)
# instead of sitting alone in the trailing paragraph at the end of the doc.

$d = $word.ActiveDocument

function Replace-WithRuns($oldText, $run1Text, $run2Text, [bool]$addGoBack) {
    $r = $d.Content
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

    $bookmarkXml = ""
    if ($addGoBack) {
        $bookmarkXml = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
    }

    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' +
           '<w:r><w:t>' + $run1Text + '</w:t></w:r>' +
           '<w:r><w:t>' + $run2Text + '</w:t></w:r>' +
           $bookmarkXml +
           '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($xml)
}

# Flight service port: 8082 -> 8080, typed as "808" + "0".
Replace-WithRuns "8082" "808" "0" $false

# Cars service port: 8083 -> 8082, typed as "808" + "2". This cell is the
# last one edited, so it now carries the "_GoBack" bookmark.
Replace-WithRuns "8083" "808" "2" $true

# Middleware port: 6659 -> 6666, typed as "66" + "66", wiring the middleware
# up to the renumbered servers.
Replace-WithRuns "6659" "66" "66" $false

# Word keeps only one "_GoBack" bookmark (the most recent edit position), so
# drop the old one that used to live by itself in the last paragraph of the
# document.
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$lastParaRange = $d.Paragraphs.Last.Range
if ($oldGoBack.Start -ge $lastParaRange.Start -and $oldGoBack.Start -le $lastParaRange.End) {
    $oldGoBack.Delete()
}
